$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Lama1"
$row2[0,2] = "Itga1"
$row2[0,3] = "ECs"
$row2[0,4] = 1
$row2[0,5] = 0.3333333333333333
$row2[0,6] = 0.014112
$row2[0,7] = 0.042336
$row2[0,8] = 0.1773673913134555
$row2[0,9] = 0.1773673913134555
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 42.105049
$row2[0,13] = 126.315147
$row2[0,14] = 0.596182887750605
$row2[0,15] = 0.5961828877506051
$row2[0,16] = 0.594186451488
$row2[0,17] = 5.347678063392
$row2[0,18] = 0.1057434035460475
$row2[0,19] = 0.1057434035460475
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object "object[,]" 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Lama1"
$row3[0,2] = "Itga1"
$row3[0,3] = "FAPs"
$row3[0,4] = 1
$row3[0,5] = 0.3333333333333333
$row3[0,6] = 0.014112
$row3[0,7] = 0.042336
$row3[0,8] = 0.1773673913134555
$row3[0,9] = 0.1773673913134555
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 2.924192
$row3[0,13] = 8.772575999999999
$row3[0,14] = 0.04140484982922635
$row3[0,15] = 0.04140484982922635
$row3[0,16] = 0.04126619750399999
$row3[0,17] = 0.3713957775359999
$row3[0,18] = 0.00734387020193525
$row3[0,19] = 0.00734387020193525
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object "object[,]" 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Lama1"
$row4[0,2] = "Itga1"
$row4[0,3] = "Inflammatory-Mac"
$row4[0,4] = 1
$row4[0,5] = 0.3333333333333333
$row4[0,6] = 0.014112
$row4[0,7] = 0.042336
$row4[0,8] = 0.1773673913134555
$row4[0,9] = 0.1773673913134555
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 7.703340666666667
$row4[0,13] = 23.110022
$row4[0,14] = 0.1090748020262369
$row4[0,15] = 0.1090748020262369
$row4[0,16] = 0.108709543488
$row4[0,17] = 0.978385891392
$row4[0,18] = 0.01934631309342525
$row4[0,19] = 0.01934631309342525
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object "object[,]" 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Lama1"
$row5[0,2] = "Itga1"
$row5[0,3] = "MuSCs"
$row5[0,4] = 1
$row5[0,5] = 0.3333333333333333
$row5[0,6] = 0.014112
$row5[0,7] = 0.042336
$row5[0,8] = 0.1773673913134555
$row5[0,9] = 0.1773673913134555
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 16.80268266666667
$row5[0,13] = 50.40804800000001
$row5[0,14] = 0.2379161671126513
$row5[0,15] = 0.2379161671126513
$row5[0,16] = 0.237119457792
$row5[0,17] = 2.134075120128
$row5[0,18] = 0.04219856991206709
$row5[0,19] = 0.04219856991206709
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object "object[,]" 1,20
$row6[0,0] = "ECs"
$row6[0,1] = "Lama1"
$row6[0,2] = "Itga1"
$row6[0,3] = "Neutrophils"
$row6[0,4] = 1
$row6[0,5] = 0.3333333333333333
$row6[0,6] = 0.014112
$row6[0,7] = 0.042336
$row6[0,8] = 0.1773673913134555
$row6[0,9] = 0.1773673913134555
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 0.6492376666666667
$row6[0,13] = 1.947713
$row6[0,14] = 0.009192825947068677
$row6[0,15] = 0.009192825947068679
$row6[0,16] = 0.009162041952
$row6[0,17] = 0.08245837756800001
$row6[0,18] = 0.001630507557030217
$row6[0,19] = 0.001630507557030217
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object "object[,]" 1,20
$row7[0,0] = "ECs"
$row7[0,1] = "Lama1"
$row7[0,2] = "Itga1"
$row7[0,3] = "Resolving-Mac"
$row7[0,4] = 1
$row7[0,5] = 0.3333333333333333
$row7[0,6] = 0.014112
$row7[0,7] = 0.042336
$row7[0,8] = 0.1773673913134555
$row7[0,9] = 0.1773673913134555
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 0.4398816666666667
$row7[0,13] = 1.319645
$row7[0,14] = 0.006228467334211686
$row7[0,15] = 0.006228467334211686
$row7[0,16] = 0.00620761008
$row7[0,17] = 0.05586849072
$row7[0,18] = 0.001104727002950199
$row7[0,19] = 0.001104727002950199
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object "object[,]" 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Lama1"
$row8[0,2] = "Itga1"
$row8[0,3] = "ECs"
$row8[0,4] = 2
$row8[0,5] = 0.6666666666666666
$row8[0,6] = 0.059263
$row8[0,7] = 0.177789
$row8[0,8] = 0.7448500362393221
$row8[0,9] = 0.7448500362393219
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 42.105049
$row8[0,13] = 126.315147
$row8[0,14] = 0.596182887750605
$row8[0,15] = 0.5961828877506051
$row8[0,16] = 2.495271518887
$row8[0,17] = 22.457443669983
$row8[0,18] = 0.4440668455463019
$row8[0,19] = 0.4440668455463018
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object "object[,]" 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Lama1"
$row9[0,2] = "Itga1"
$row9[0,3] = "FAPs"
$row9[0,4] = 2
$row9[0,5] = 0.6666666666666666
$row9[0,6] = 0.059263
$row9[0,7] = 0.177789
$row9[0,8] = 0.7448500362393221
$row9[0,9] = 0.7448500362393219
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 2.924192
$row9[0,13] = 8.772575999999999
$row9[0,14] = 0.04140484982922635
$row9[0,15] = 0.04140484982922635
$row9[0,16] = 0.173296390496
$row9[0,17] = 1.559667514464
$row9[0,18] = 0.03084040389578294
$row9[0,19] = 0.03084040389578293
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object "object[,]" 1,20
$row10[0,0] = "FAPs"
$row10[0,1] = "Lama1"
$row10[0,2] = "Itga1"
$row10[0,3] = "Inflammatory-Mac"
$row10[0,4] = 2
$row10[0,5] = 0.6666666666666666
$row10[0,6] = 0.059263
$row10[0,7] = 0.177789
$row10[0,8] = 0.7448500362393221
$row10[0,9] = 0.7448500362393219
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 7.703340666666667
$row10[0,13] = 23.110022
$row10[0,14] = 0.1090748020262369
$row10[0,15] = 0.1090748020262369
$row10[0,16] = 0.4565230779286667
$row10[0,17] = 4.108707701358
$row10[0,18] = 0.08124437024203944
$row10[0,19] = 0.08124437024203943
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object "object[,]" 1,20
$row11[0,0] = "FAPs"
$row11[0,1] = "Lama1"
$row11[0,2] = "Itga1"
$row11[0,3] = "MuSCs"
$row11[0,4] = 2
$row11[0,5] = 0.6666666666666666
$row11[0,6] = 0.059263
$row11[0,7] = 0.177789
$row11[0,8] = 0.7448500362393221
$row11[0,9] = 0.7448500362393219
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 16.80268266666667
$row11[0,13] = 50.40804800000001
$row11[0,14] = 0.2379161671126513
$row11[0,15] = 0.2379161671126513
$row11[0,16] = 0.9957773828746669
$row11[0,17] = 8.961996445872002
$row11[0,18] = 0.1772118656957789
$row11[0,19] = 0.1772118656957789
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object "object[,]" 1,20
$row12[0,0] = "FAPs"
$row12[0,1] = "Lama1"
$row12[0,2] = "Itga1"
$row12[0,3] = "Neutrophils"
$row12[0,4] = 2
$row12[0,5] = 0.6666666666666666
$row12[0,6] = 0.059263
$row12[0,7] = 0.177789
$row12[0,8] = 0.7448500362393221
$row12[0,9] = 0.7448500362393219
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 0.6492376666666667
$row12[0,13] = 1.947713
$row12[0,14] = 0.009192825947068677
$row12[0,15] = 0.009192825947068679
$row12[0,16] = 0.03847577183966667
$row12[0,17] = 0.3462819465570001
$row12[0,18] = 0.006847276739815885
$row12[0,19] = 0.006847276739815884
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object "object[,]" 1,20
$row13[0,0] = "FAPs"
$row13[0,1] = "Lama1"
$row13[0,2] = "Itga1"
$row13[0,3] = "Resolving-Mac"
$row13[0,4] = 2
$row13[0,5] = 0.6666666666666666
$row13[0,6] = 0.059263
$row13[0,7] = 0.177789
$row13[0,8] = 0.7448500362393221
$row13[0,9] = 0.7448500362393219
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 0.4398816666666667
$row13[0,13] = 1.319645
$row13[0,14] = 0.006228467334211686
$row13[0,15] = 0.006228467334211686
$row13[0,16] = 0.02606870721166667
$row13[0,17] = 0.234618364905
$row13[0,18] = 0.004639274119603008
$row13[0,19] = 0.004639274119603007
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object "object[,]" 1,20
$row14[0,0] = "MuSCs"
$row14[0,1] = "Lama1"
$row14[0,2] = "Itga1"
$row14[0,3] = "ECs"
$row14[0,4] = 1
$row14[0,5] = 0.3333333333333333
$row14[0,6] = 0.006188666666666666
$row14[0,7] = 0.018566
$row14[0,8] = 0.07778257244722256
$row14[0,9] = 0.07778257244722254
$row14[0,10] = 3
$row14[0,11] = 1
$row14[0,12] = 42.105049
$row14[0,13] = 126.315147
$row14[0,14] = 0.596182887750605
$row14[0,15] = 0.5961828877506051
$row14[0,16] = 0.2605741132446667
$row14[0,17] = 2.345167019202
$row14[0,18] = 0.04637263865825579
$row14[0,19] = 0.04637263865825579
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object "object[,]" 1,20
$row15[0,0] = "MuSCs"
$row15[0,1] = "Lama1"
$row15[0,2] = "Itga1"
$row15[0,3] = "FAPs"
$row15[0,4] = 1
$row15[0,5] = 0.3333333333333333
$row15[0,6] = 0.006188666666666666
$row15[0,7] = 0.018566
$row15[0,8] = 0.07778257244722256
$row15[0,9] = 0.07778257244722254
$row15[0,10] = 3
$row15[0,11] = 1
$row15[0,12] = 2.924192
$row15[0,13] = 8.772575999999999
$row15[0,14] = 0.04140484982922635
$row15[0,15] = 0.04140484982922635
$row15[0,16] = 0.01809684955733333
$row15[0,17] = 0.162871646016
$row15[0,18] = 0.003220575731508169
$row15[0,19] = 0.003220575731508169
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object "object[,]" 1,20
$row16[0,0] = "MuSCs"
$row16[0,1] = "Lama1"
$row16[0,2] = "Itga1"
$row16[0,3] = "Inflammatory-Mac"
$row16[0,4] = 1
$row16[0,5] = 0.3333333333333333
$row16[0,6] = 0.006188666666666666
$row16[0,7] = 0.018566
$row16[0,8] = 0.07778257244722256
$row16[0,9] = 0.07778257244722254
$row16[0,10] = 3
$row16[0,11] = 1
$row16[0,12] = 7.703340666666667
$row16[0,13] = 23.110022
$row16[0,14] = 0.1090748020262369
$row16[0,15] = 0.1090748020262369
$row16[0,16] = 0.04767340760577778
$row16[0,17] = 0.429060668452
$row16[0,18] = 0.00848411869077223
$row16[0,19] = 0.00848411869077223
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object "object[,]" 1,20
$row17[0,0] = "MuSCs"
$row17[0,1] = "Lama1"
$row17[0,2] = "Itga1"
$row17[0,3] = "MuSCs"
$row17[0,4] = 1
$row17[0,5] = 0.3333333333333333
$row17[0,6] = 0.006188666666666666
$row17[0,7] = 0.018566
$row17[0,8] = 0.07778257244722256
$row17[0,9] = 0.07778257244722254
$row17[0,10] = 3
$row17[0,11] = 1
$row17[0,12] = 16.80268266666667
$row17[0,13] = 50.40804800000001
$row17[0,14] = 0.2379161671126513
$row17[0,15] = 0.2379161671126513
$row17[0,16] = 0.1039862021297778
$row17[0,17] = 0.9358758191680001
$row17[0,18] = 0.01850573150480531
$row17[0,19] = 0.01850573150480531
$ws.Range("A17:T17").Value = $row17

$row18 = New-Object "object[,]" 1,20
$row18[0,0] = "MuSCs"
$row18[0,1] = "Lama1"
$row18[0,2] = "Itga1"
$row18[0,3] = "Neutrophils"
$row18[0,4] = 1
$row18[0,5] = 0.3333333333333333
$row18[0,6] = 0.006188666666666666
$row18[0,7] = 0.018566
$row18[0,8] = 0.07778257244722256
$row18[0,9] = 0.07778257244722254
$row18[0,10] = 3
$row18[0,11] = 1
$row18[0,12] = 0.6492376666666667
$row18[0,13] = 1.947713
$row18[0,14] = 0.009192825947068677
$row18[0,15] = 0.009192825947068679
$row18[0,16] = 0.004017915506444444
$row18[0,17] = 0.03616123955800001
$row18[0,18] = 0.0007150416502225767
$row18[0,19] = 0.0007150416502225767
$ws.Range("A18:T18").Value = $row18

$row19 = New-Object "object[,]" 1,20
$row19[0,0] = "MuSCs"
$row19[0,1] = "Lama1"
$row19[0,2] = "Itga1"
$row19[0,3] = "Resolving-Mac"
$row19[0,4] = 1
$row19[0,5] = 0.3333333333333333
$row19[0,6] = 0.006188666666666666
$row19[0,7] = 0.018566
$row19[0,8] = 0.07778257244722256
$row19[0,9] = 0.07778257244722254
$row19[0,10] = 3
$row19[0,11] = 1
$row19[0,12] = 0.4398816666666667
$row19[0,13] = 1.319645
$row19[0,14] = 0.006228467334211686
$row19[0,15] = 0.006228467334211686
$row19[0,16] = 0.002722281007777778
$row19[0,17] = 0.02450052907
$row19[0,18] = 0.0004844662116584796
$row19[0,19] = 0.0004844662116584795
$ws.Range("A19:T19").Value = $row19

Write-Output "done"